$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = [double]99519589
$ws.Range("B3").Value = [double]93235
$ws.Range("E3").Value = [double]210
$ws.Range("F3").Value = 'Grön sköldmossa'
$ws.Range("G3").Value = 'Buxbaumia viridis'
$ws.Range("H3").Value = '(Moug. ex Lam. & DC.) Brid. ex Moug. & Nestl.'
$ws.Range("K3").ClearContents() | Out-Null
$ws.Range("L3").ClearContents() | Out-Null
$ws.Range("P3").Value = 'Norr om Varnö, Srm'
$ws.Range("Q3").Value = [double]693062.2900751554
$ws.Range("R3").Value = [double]6549110.554353707
$ws.Range("S3").Value = [double]20
$ws.Range("Z3").Value = '00:00'
$ws.Range("AB3").Value = '00:00'
$ws.Range("AC3").Value = '3 exemplar'
$ws.Range("AW3").Value = 'Amanda Tas'
$ws.Range("AX3").Value = 'Amanda Tas'

# Row 4
$ws.Range("A4").Value = [double]99519598
$ws.Range("B4").Value = [double]98520
$ws.Range("E4").Value = [double]222498
$ws.Range("F4").Value = 'Blåsippa'
$ws.Range("G4").Value = 'Hepatica nobilis'
$ws.Range("H4").Value = 'Schreb.'
$ws.Range("Q4").Value = [double]693046.2694973514
$ws.Range("R4").Value = [double]6549110.77766049
$ws.Range("AC4").ClearContents() | Out-Null

# Row 5
$ws.Range("A5").Value = [double]99504553
$ws.Range("B5").Value = [double]57064
$ws.Range("D5").Value = 'NT'
$ws.Range("E5").Value = [double]103055
$ws.Range("F5").Value = 'Gulsparv'
$ws.Range("G5").Value = 'Emberiza citrinella'
$ws.Range("H5").Value = 'Linnaeus, 1758'
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = '1'
$ws.Range("J5").ClearContents() | Out-Null
$ws.Range("K5").Value = 'adult'
$ws.Range("M5").Value = 'spel/sång'
$ws.Range("P5").Value = 'Risselrum, Ornö, Srm'
$ws.Range("Q5").Value = [double]693113.6323470459
$ws.Range("R5").Value = [double]6549157.477006816
$ws.Range("S5").Value = [double]10
$ws.Range("Z5").Value = '14:30'
$ws.Range("AB5").Value = '14:30'
$ws.Range("AC5").Value = 'Sjungande gulsparv'
$ws.Range("AF5").ClearContents() | Out-Null
$ws.Range("AI5").Value = 'Blandskog'

# Row 6
$ws.Range("A6").Value = [double]99505564
$ws.Range("B6").Value = [double]56411
$ws.Range("D6").Value = 'NT'
$ws.Range("E6").Value = [double]100049
$ws.Range("F6").Value = 'Spillkråka'
$ws.Range("G6").Value = 'Dryocopus martius'
$ws.Range("H6").Value = '(Linnaeus, 1758)'
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = '1'
$ws.Range("K6").Value = ""
$ws.Range("L6").Value = ""
$ws.Range("M6").Value = 'födosökande'
$ws.Range("N6").Value = ""
$ws.Range("P6").Value = 'Risselrum, Ornö, Srm'
$ws.Range("Q6").Value = [double]693074.6979561135
$ws.Range("R6").Value = [double]6549130.769313719
$ws.Range("S6").Value = [double]10
$ws.Range("Z6").Value = '14:00'
$ws.Range("AB6").Value = '14:00'
$ws.Range("AC6").Value = 'Spillkråka observerades vid granlåga med färska hackspår'
$ws.Range("AH6").Value = 'Skogsmark'
$ws.Range("AI6").Value = 'Blandskog'
$ws.Range("AW6").Value = 'Michael Lander'
$ws.Range("AX6").Value = 'Michael Lander'

# Row 7
$ws.Range("A7").Value = [double]99519597
$ws.Range("Q7").Value = [double]693062.5255223531
$ws.Range("R7").Value = [double]6549177.581972098

# Row 8
$ws.Range("A8").Value = [double]99519606
$ws.Range("B8").Value = [double]103813
$ws.Range("D8").Value = 'EN'
$ws.Range("E8").Value = [double]220785
$ws.Range("F8").Value = 'Ask'
$ws.Range("G8").Value = 'Fraxinus excelsior'
$ws.Range("H8").Value = 'L.'
$ws.Range("Q8").Value = [double]693083.1766483777
$ws.Range("R8").Value = [double]6549116.247088743
$ws.Range("AC8").Value = 'Föryngring'

# Row 9
$ws.Range("A9").Value = [double]99486685
$ws.Range("B9").Value = [double]98520
$ws.Range("D9").Value = 'LC'
$ws.Range("E9").Value = [double]222498
$ws.Range("F9").Value = 'Blåsippa'
$ws.Range("G9").Value = 'Hepatica nobilis'
$ws.Range("H9").Value = 'Schreb.'
$ws.Range("I9").Value = ""
$ws.Range("J9").ClearContents() | Out-Null
$ws.Range("L9").Value = ""
$ws.Range("N9").ClearContents() | Out-Null
$ws.Range("P9").Value = 'Norr om Varnö, Ornö, Srm'
$ws.Range("Q9").Value = [double]692892.0329120732
$ws.Range("R9").Value = [double]6549045.783047367
$ws.Range("S9").Value = [double]10
$ws.Range("Z9").Value = '14:07'
$ws.Range("AB9").Value = '14:07'
$ws.Range("AF9").ClearContents() | Out-Null
$ws.Range("AH9").ClearContents() | Out-Null
$ws.Range("AI9").ClearContents() | Out-Null
$ws.Range("AJ9").ClearContents() | Out-Null
$ws.Range("AK9").ClearContents() | Out-Null
$ws.Range("AL9").ClearContents() | Out-Null
$ws.Range("AO9").ClearContents() | Out-Null
$ws.Range("AW9").Value = 'Kristina Bäck'
$ws.Range("AX9").Value = 'Kristina Bäck'

# Row 10
$ws.Range("A10").Value = [double]99505966
$ws.Range("B10").Value = [double]98520
$ws.Range("D10").Value = 'LC'
$ws.Range("E10").Value = [double]222498
$ws.Range("F10").Value = 'Blåsippa'
$ws.Range("G10").Value = 'Hepatica nobilis'
$ws.Range("H10").Value = 'Schreb.'
$ws.Range("I10").Value = ""
$ws.Range("J10").Value = ""
$ws.Range("K10").Value = 'blomning'
$ws.Range("M10").ClearContents() | Out-Null
$ws.Range("P10").Value = 'Storrum, Ornö, Srm'
$ws.Range("Q10").Value = [double]692884.3662214879
$ws.Range("R10").Value = [double]6549105.711763832
$ws.Range("S10").Value = [double]5
$ws.Range("Z10").Value = '13:00'
$ws.Range("AB10").Value = '13:00'
$ws.Range("AC10").Value = 'Blåsippor i ett område om ca 20 m omkrets'
$ws.Range("AF10").Value = ""
$ws.Range("AI10").Value = 'Blandskog med inslag av gamla ekar'

# Row 11
$ws.Range("A11").Value = [double]99519600
$ws.Range("B11").Value = [double]98520
$ws.Range("D11").Value = 'LC'
$ws.Range("E11").Value = [double]222498
$ws.Range("F11").Value = 'Blåsippa'
$ws.Range("G11").Value = 'Hepatica nobilis'
$ws.Range("H11").Value = 'Schreb.'
$ws.Range("Q11").Value = [double]692872.8119970543
$ws.Range("R11").Value = [double]6549068.528527547
$ws.Range("AC11").ClearContents() | Out-Null

# Row 12
$ws.Range("A12").Value = [double]99519601
$ws.Range("B12").Value = [double]98520
$ws.Range("D12").Value = 'LC'
$ws.Range("E12").Value = [double]222498
$ws.Range("F12").Value = 'Blåsippa'
$ws.Range("G12").Value = 'Hepatica nobilis'
$ws.Range("H12").Value = 'Schreb.'
$ws.Range("I12").Value = ""
$ws.Range("K12").ClearContents() | Out-Null
$ws.Range("L12").ClearContents() | Out-Null
$ws.Range("M12").ClearContents() | Out-Null
$ws.Range("N12").ClearContents() | Out-Null
$ws.Range("P12").Value = 'Norr om Varnö, Srm'
$ws.Range("Q12").Value = [double]692865.2778250941
$ws.Range("R12").Value = [double]6549023.299178503
$ws.Range("S12").Value = [double]20
$ws.Range("Z12").Value = '00:00'
$ws.Range("AB12").Value = '00:00'
$ws.Range("AC12").ClearContents() | Out-Null
$ws.Range("AH12").ClearContents() | Out-Null
$ws.Range("AI12").ClearContents() | Out-Null
$ws.Range("AW12").Value = 'Amanda Tas'
$ws.Range("AX12").Value = 'Amanda Tas'

# Row 13
$ws.Range("A13").Value = [double]99505991
$ws.Range("B13").Value = [double]89412
$ws.Range("D13").Value = 'NT'
$ws.Range("E13").Value = [double]5442
$ws.Range("F13").Value = 'Tallticka'
$ws.Range("G13").Value = 'Porodaedalea pini'
$ws.Range("H13").Value = '(Brot.) Murrill'
$ws.Range("I13").NumberFormat = "@"
$ws.Range("I13").Value = '1'
$ws.Range("J13").Value = ""
$ws.Range("K13").Value = ""
$ws.Range("N13").Value = ""
$ws.Range("P13").Value = 'Storrum, Ornö, Srm'
$ws.Range("Q13").Value = [double]692873.4220747473
$ws.Range("R13").Value = [double]6549097.427865601
$ws.Range("S13").Value = [double]5
$ws.Range("Z13").Value = '13:30'
$ws.Range("AB13").Value = '14:00'
$ws.Range("AF13").Value = ""
$ws.Range("AH13").Value = 'Skogsmark'
$ws.Range("AI13").Value = 'Blandskog'
$ws.Range("AJ13").Value = 'tall'
$ws.Range("AK13").Value = 'Pinus sylvestris'
$ws.Range("AL13").Value = 'Växer på död tall'
$ws.Range("AO13").Value = 'Pinus sylvestris # Växer på död tall'
$ws.Range("AW13").Value = 'Michael Lander'
$ws.Range("AX13").Value = 'Michael Lander'

# Row 14
$ws.Range("A14").Value = [double]99519607
$ws.Range("Q14").Value = [double]692885.6932734415
$ws.Range("R14").Value = [double]6548966.588938461

